# Re-map the "Address" column (D) to the refreshed mock-data set, resize the
# "Unit ID" column (C), and move the active selection - per the
# "remake json feature" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Address (column D) values for each data row, replacing the old
# street/city/province/postal-code strings with the regenerated ones.
$addresses = @{
    2  = "9998 Birch Blvd., Toronto, NL, L5M6N4"
    3  = "563 Queen St., St. John's, BC, M6H 9J0"
    4  = "3435 Cedar Ln., St. John's, AB, P0O9K8"
    5  = "766 Birch Blvd., charletown PE, L8K7J6"
    6  = "1654 Queen St., Halifax, NB, P1K3L4"
    7  = "4947 Maple St., Regina , SK"
    8  = "119 Birch Blvd., Calgary, SK"
    9  = "4899 Queen St., Vancouver"
    10 = "134 Pine Ave., Winnipeg, MB, "
    11 = "7752 King Rd., Winnipeg,"
    12 = "3130 Elm Dr., Fredericton, "
    13 = "4278 Elm Dr., Charlottetown, BC"
    14 = "579 King Rd., Regina, "
    15 = "3200 Cedar Ln., Calgary,"
    16 = "4483 Oak St., Fredericton, BC"
    17 = "9801 Elm Dr., Fredericton, NL"
    18 = "6029 Main St., Vancouver,"
    19 = "2789 Oak St., Calgary"
    20 = "4258 Oak St., Toronto, PE"
    21 = "1685 Queen St., Charlottetown,, R9S 4C1"
}

foreach ($row in $addresses.Keys) {
    $ws.Cells.Item($row, 4).Value = $addresses[$row]
}

# Widen the "Unit ID" column (C) to fit its contents.
$ws.Columns.Item(3).ColumnWidth = 22.71

# Move the active selection from D16 to I5.
$ws.Range("I5").Select() | Out-Null
